$d = $word.ActiveDocument

# Locate the placeholder room text "Room: TBA." and split it out so the
# final result reads "Room: 5125 (Edison) - 408E." instead.
$roomText = "5125 (Edison) - 408E"

$full = $d.Content
$found = $full.Find.Execute("Room: TBA.")
if (-not $found) {
    throw "Could not find placeholder text 'Room: TBA.'"
}
$paraStart = $full.Start
$paraEnd = $full.End

# Sub-ranges for the fixed label, the "TBA" placeholder, and the trailing
# period, all computed relative to the match (not hardcoded offsets).
$labelStart = $paraStart
$labelEnd = $paraStart + 6          # "Room: " (6 chars)
$placeholderStart = $labelEnd
$placeholderEnd = $paraEnd - 1      # up to, but not including, the "."
$periodStart = $placeholderEnd
$periodEnd = $paraEnd

# The _GoBack bookmark currently sits right after the trailing period; move
# it so it instead sits right after the placeholder (i.e. right where the
# new room value will end), matching where Word leaves it after the user
# types the replacement text there.
$bmTarget = $d.Range($placeholderEnd, $placeholderEnd)
$d.Bookmarks.Add("_GoBack", $bmTarget)

# Replace the "TBA" placeholder text with the real room value.
$placeholder = $d.Range($placeholderStart, $placeholderEnd)
$placeholder.Text = $roomText

# Split the merged run into distinct runs for "Room: " / room value / "."
# by toggling a character property off again (net-zero formatting change,
# but it forces the engine to keep these as separate <w:r> runs instead of
# silently re-merging them).
$label = $d.Range($labelStart, $labelEnd)
$label.Font.Bold = 1
$label.Font.Bold = 0

$value = $d.Range($placeholderStart, $placeholderStart + $roomText.Length)
$value.Font.Bold = 1
$value.Font.Bold = 0

Write-Output "Room text updated to: $roomText"
